$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 21
$ws.Range("C4").Value = 78
$ws.Range("B4").Value = "26/01/2024"

$ws.Range("C3:C4").NumberFormat = "0.00"
$ws.Range("C3").Font.Underline = $true

$ws.Range("C3").Select()
